$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the conversion text in A1 ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$cellA1 = $wsHoja1.Range("A1")
$oldText = $cellA1.Value2
$newText = $oldText -replace "1000 Bs = 5.41 = 21559.46 pesos", "1000 Bs = 5.38 = 21406.62 pesos"
$newText = $newText -replace "21559.46 pesos = 5.39 = 969.61 Bs", "21406.62 pesos = 5.34 = 945.41 Bs"
$cellA1.Value2 = $newText

# --- Sheet "tasas": update rate values ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value2 = 185.758
$wsTasas.Range("O10").Value2 = 3976.45
$wsTasas.Range("N12").Value2 = 4010
$wsTasas.Range("O12").Value2 = 177.1
